$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = -12.3182
$ws.Range("C3").Value = -11.9211
$ws.Range("C5").Value = -12.0217
$ws.Range("D7").Value = -7.313599999999987
$ws.Range("A9").Value = -20.28819999999997
$ws.Range("D9").Value = -8.811700000000004
$ws.Range("C11").Value = -14.10300000000001
$ws.Range("C12").Value = -13.8107
$ws.Range("A13").Value = -22.01990000000002
$ws.Range("A16").Value = -20.20219999999999
$ws.Range("A18").Value = -21.97170000000002
$ws.Range("A20").Value = -22.16800000000002
$ws.Range("C21").Value = -11.7683
$ws.Range("D21").Value = -7.506300000000002
